# Update the cost figures on Sheet1 with the refreshed per-project totals
# (data now saved/consolidated from a single approval.xlsx source, so the
# underlying numbers changed on this save).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4117
$ws.Range("B2").Value = 2212
$ws.Range("C2").Value = 3796
$ws.Range("D2").Value = 4830
$ws.Range("E2").Value = 2256
$ws.Range("F2").Value = 2149
